$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to fill into columns D (sex = "U"), F (sire), G (dam) for rows 2-11.
# E column (environ/year) already has data.
$rows = @(
    @{ Row = 2;  Sire = 1;  Dam = 1 },
    @{ Row = 3;  Sire = 3;  Dam = 1 },
    @{ Row = 4;  Sire = 5;  Dam = 1 },
    @{ Row = 5;  Sire = 3;  Dam = 5 },
    @{ Row = 6;  Sire = 5;  Dam = 5 },
    @{ Row = 7;  Sire = 6;  Dam = 6 },
    @{ Row = 8;  Sire = 8;  Dam = 6 },
    @{ Row = 9;  Sire = 10; Dam = 6 },
    @{ Row = 10; Sire = 8;  Dam = 10 },
    @{ Row = 11; Sire = 10; Dam = 10 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = "U"
    $ws.Cells.Item($r.Row, 6).Value = $r.Sire
    $ws.Cells.Item($r.Row, 7).Value = $r.Dam
}

# Update the selection to reflect the new active range
$ws.Range("D2:G11").Select()
